# Updates Price (D) and Volume(1h) (E) columns for the cryptos worksheet.
# Values are written as text (matching the original inlineStr cell type); numeric-
# looking prices are written with a leading apostrophe and then the cell style is
# reset to "Normal" so Excel stores them as text without altering cell formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.624.22"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "1.641.63"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.28%  "

$cell = $ws.Range("D5")
$cell.Value = "'215.11"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +0.63%  "

$cell = $ws.Range("D10")
$cell.Value = "'19.24"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "1.871.01"
$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").Value = "1.636.05"
$ws.Range("E14").Value = "  +0.74%  "

$cell = $ws.Range("D15")
$cell.Value = "'0.529"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "

$cell = $ws.Range("D16")
$cell.Value = "'65.37"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.87%  "

$ws.Range("D17").Value = "26.666.10"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "0.0₃0746"
$ws.Range("E18").Value = "  +0.70%  "

$cell = $ws.Range("D19")
$cell.Value = "'216.37"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("E20").Value = "  +0.23%  "

$cell = $ws.Range("D21")
$cell.Value = "'4.35"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "

$cell = $ws.Range("D22")
$cell.Value = "'6.29"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.13%  "

$cell = $ws.Range("D23")
$cell.Value = "'9.51"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "

$cell = $ws.Range("D24")
$cell.Value = "'2.20"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +12.75%  "

$cell = $ws.Range("D25")
$cell.Value = "'145.55"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  -0.86%  "

$cell = $ws.Range("D28")
$cell.Value = "'7.15"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +4.21%  "

$ws.Range("E29").Value = "  +1.51%  "

$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("E31").Value = "  +0.30%  "

$cell = $ws.Range("D32")
$cell.Value = "'3.38"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("D34").Value = "1.278.40"
$ws.Range("E34").Value = "  +4.65%  "

$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("E36").Value = "  +5.64%  "

$ws.Range("E37").Value = "  +0.39%  "

$cell = $ws.Range("D38")
$cell.Value = "'0.531"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +6.11%  "

$cell = $ws.Range("D39")
$cell.Value = "'0.827"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.46%  "

$cell = $ws.Range("D40")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "

$cell = $ws.Range("D41")
$cell.Value = "'0.817"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +2.71%  "

$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("E43").Value = "  +2.48%  "

$ws.Range("D44").Value = "1.781.17"
$ws.Range("E44").Value = "  +0.85%  "

$cell = $ws.Range("D45")
$cell.Value = "'92.02"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "

$cell = $ws.Range("D46")
$cell.Value = "'59.91"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +8.29%  "

$cell = $ws.Range("D47")
$cell.Value = "'1.59"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("E48").Value = "  +0.87%  "

$cell = $ws.Range("D49")
$cell.Value = "'7.77"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.52%  "

$cell = $ws.Range("D50")
$cell.Value = "'0.0969"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.79%  "

$ws.Range("E51").Value = "  -0.28%  "
